# The external data refresh (sector performance add-in) reordered rows
# 2-12 and moved the lone annotated cell (column G) from row 2 to row 8.
# D/E columns are untouched; A/B/C are rewritten per the new row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newOrder = @(
    @{ Row = 2;  Stock = "^SPTTEN";  B = 267.13;              C = 267.13 },
    @{ Row = 3;  Stock = "^GSPTTMT"; B = 378.22;              C = 378.22 },
    @{ Row = 4;  Stock = "^GSPTTIN"; B = 458.44;              C = 458.44 },
    @{ Row = 5;  Stock = "^GSPTTCD"; B = 298.46;              C = 298.46 },
    @{ Row = 6;  Stock = "^GSPTTCS"; B = 1062.37;             C = 1062.37 },
    @{ Row = 7;  Stock = "^GSPTTHC"; B = 25.73;               C = 25.73 },
    @{ Row = 8;  Stock = "^SPTTFS";  B = 478.51281531759696;  C = 472.34 },
    @{ Row = 9;  Stock = "^SPTTTK";  B = 281.82;              C = 281.82 },
    @{ Row = 10; Stock = "^GSPTTTS"; B = 137.64;              C = 137.64 },
    @{ Row = 11; Stock = "^GSPTTUT"; B = 306.67;              C = 804.45 },
    @{ Row = 12; Stock = "^GSPRTRE"; B = 314.57;              C = 314.57 }
)

foreach ($item in $newOrder) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.Stock
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
}

# Rows 3-12 (the non-"boy-of-year" rows) lose their fill-applying style
# variant (s=8) in favour of the plain centred style (s=1) already used
# elsewhere - clearing the (no-op) interior pattern reassigns them without
# introducing a new style entry.
$ws.Range("B3:C12").Interior.Pattern = -4142

# The annotated/empty helper cell in column G moves from row 2 to row 8.
$ws.Range("G2").Copy($ws.Range("G8"))
$ws.Range("G2").ClearFormats()

Write-Host "edit applied"
